# Add Function Points into UseCaseDoc
# - Convert the H3:H7 "Total" column from text ("27h", "2h", ...) to plain
#   numeric hour values and fill in the previously-empty "FP" (I) column
#   with the computed function-point numbers.
# - Highlight the still-unfilled placeholder cells (H9:I11) in red so they
#   stand out as TODO.
# - Fill in the Function Point total for UC9 (I12).
# - Add a scatter chart (with linear trendline) plotting FP against Total
#   hours, sourced from H3:I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Total (H) / FP (I) values for UC1..UC5 -------------------------------
$ws.Range("H3").Value = 27
$ws.Range("I3").Value = 102.6

$ws.Range("H4").Value = 27
$ws.Range("I4").Value = 50.92

$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 40.28

$ws.Range("H6").Value = 7.5
$ws.Range("I6").Value = 34.96

$ws.Range("H7").Value = 2.5
$ws.Range("I7").Value = 24.32

# --- Function Point total for UC9 ------------------------------------------
$ws.Range("I12").Value = 20.52

# --- Mark still-empty placeholder cells in red so they are easy to spot ----
$ws.Range("H9").Font.Color = 255
$ws.Range("I9").Font.Color = 255
$ws.Range("H10").Font.Color = 255
$ws.Range("I10").Font.Color = 255
$ws.Range("H11").Font.Color = 255
$ws.Range("I11").Font.Color = 255

# --- Scatter chart: FP (y) vs Total hours (x), UC1..UC5, with trendline ----
$chartObj = $ws.Shapes.AddChart2(-1, 74)
$chart = $chartObj.Chart
$chart.SetSourceData($ws.Range("H3:I7"))

$series = $chart.SeriesCollection(1)
$series.Name = '=Tabelle1!$I$2'
$series.Format.Line.Weight = 2.25

$trendlines = $series.Trendlines()
$trendline = $trendlines.Add()
$trendline.DisplayEquation = $false
$trendline.DisplayRSquared = $false

$chart.HasLegend = $true
$chart.Legend.Position = -4152
